# Apply the "reran simulation" edit:
#  - Insert two new HKL rows ("Holden", "Rizzie Spiral") right after the
#    "Spiral5" row (i.e. before the old row 4), pushing every row below it
#    down by two.
#  - Fill in freshly-computed C:T values for the two new rows.
#  - Rename the "Thomas Hex" entry to "Matthies Hex" (same row, now shifted
#    down two rows by the insert above).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert two blank rows at row 4 (shifts old rows 4:29 -> 6:31) ----
$ws.Rows("4:5").Insert()

# --- 2. Populate the two new rows -----------------------------------------
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "Holden"
$ws.Range("A4").Font.Bold = $true
$ws.Range("B4").Font.Bold = $true

$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "Rizzie Spiral"
$ws.Range("A5").Font.Bold = $true
$ws.Range("B5").Font.Bold = $true

$row4vals = @(0.9822041070643793,1.003803562701073,1.005610714995183,0.9822041070643793,1.009798607245047,0.9906304877670999,1.005610714995183,1.004330426293581,1.005610714995183,1.003803562701073,0.9930038348827261,0.9930038348827261,0.9922127191775174,0.9972061282535449,0.9972061282535449,0.9993072749389544,0.9993072749389544,0.9993963176777272)
$row5vals = @(0.9700504114319629,1.009913534414464,1.003120332843378,0.9700504114319629,1.027730049476684,0.9819555917340645,1.003120332843378,1.007933013049084,1.003120332843378,1.009913534414464,0.9899819729232135,0.9899819729232135,0.9873065125268305,0.994361426229935,0.994361426229935,0.9965511528832957,0.9965511528832957,1.000117155491606)

$cols = @("C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "4").Value = $row4vals[$i]
    $ws.Range($cols[$i] + "5").Value = $row5vals[$i]
}

# --- 3. Rename "Thomas Hex" -> "Matthies Hex" (now two rows further down) -
$found = $ws.Cells.Find("Thomas Hex")
if ($found -ne $null) {
    $found.Value = "Matthies Hex"
} else {
    $ws.Range("B11").Value = "Matthies Hex"
}
